# Update "想去人数" (interested-count) figures across the three sheets that
# track exhibitions (展览), performances (演出), and the combined listing
# (全部类型). Only column F values change; everything else stays the same.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 908
$ws1.Range("F7").Value = 907
$ws1.Range("F8").Value = 711
$ws1.Range("F9").Value = 163
$ws1.Range("F11").Value = 76
$ws1.Range("F12").Value = 756
$ws1.Range("F15").Value = 482
$ws1.Range("F16").Value = 1272
$ws1.Range("F18").Value = 364
$ws1.Range("F19").Value = 1025
$ws1.Range("F20").Value = 2763
$ws1.Range("F21").Value = 1236
$ws1.Range("F22").Value = 630
$ws1.Range("F23").Value = 162
$ws1.Range("F24").Value = 1228
$ws1.Range("F26").Value = 953
$ws1.Range("F27").Value = 309
$ws1.Range("F28").Value = 159
$ws1.Range("F29").Value = 1289

# --- Sheet "演出" (performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 508
$ws2.Range("F5").Value = 8

# --- Sheet "全部类型" (combined listing) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 508
$ws4.Range("F8").Value = 508
$ws4.Range("F10").Value = 8
$ws4.Range("F12").Value = 908
$ws4.Range("F15").Value = 907
$ws4.Range("F16").Value = 711
$ws4.Range("F17").Value = 163
$ws4.Range("F23").Value = 76
$ws4.Range("F25").Value = 756
$ws4.Range("F28").Value = 482
$ws4.Range("F29").Value = 1272
$ws4.Range("F31").Value = 364
$ws4.Range("F32").Value = 1025
$ws4.Range("F33").Value = 2763
$ws4.Range("F34").Value = 1236
$ws4.Range("F35").Value = 630
$ws4.Range("F36").Value = 162
$ws4.Range("F37").Value = 1228
$ws4.Range("F40").Value = 953
$ws4.Range("F41").Value = 309
$ws4.Range("F42").Value = 159
$ws4.Range("F43").Value = 1289
